# profile ENW to Neon onboarding scripts implementation
# Adds three new ENW-header-profile-dropdown test cases (Profile64, Profile65,
# Profile66) as rows 65-67 of the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing row's formatting (borders/fonts) down onto the new rows
# before filling in values, same as the last data row in the sheet.
$ws.Range("A64:E64").Copy()
$ws.Range("A65:E67").PasteSpecial(-4122)

# Row 65: Feedback page link test case
$ws.Range("A65").Value = "Profile64"
$ws.Range("B65").Value = "OPQA-1722"
$ws.Range("C65").Value = "Verify that  A signed in ENW user shall be able to try to access their Feedback page via a link within the ENW header profile dropdown"
$ws.Range("D65").Value = "Y"

# TCIDs for rows 66 and 67 entered together
$ws.Range("A66").Value = "Profile65"
$ws.Range("A67").Value = "Profile66"

# Row 66: Privacy page link test case
$ws.Range("B66").Value = "OPQA-1715"
$ws.Range("C66").Value = "Verify that  A signed in ENW user shall be able to try to access their Privacy page via a link within the ENW header profile dropdown."
$ws.Range("D66").Value = "Y"

# Row 67: Acceptable Use page link test case
$ws.Range("B67").Value = "OPQA-1718"
$ws.Range("C67").Value = "Verify that  A signed in ENW user shall be able to try to access their Acceptable Use page via a link within the ENW header profile dropdown"
$ws.Range("D67").Value = "Y"

# Leave the cursor where the author's did when they finished editing.
[void]$ws.Range("D67").Select()
